# Fruta / hortaliza, semanal
#
# Two new weekly price-report rows for Naranja / Navel Late (Vega Modelo de
# Temuco) need to be inserted right before the current row 469, pushing the
# existing data (old rows 469-529) down to rows 471-531.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 469 (shifts everything below down by one,
# each call).
$ws.Rows.Item(469).Insert()
$ws.Rows.Item(469).Insert()

# New row 469: Navel Late / Primera, $/caja 15 kilos granel
$ws.Range("A469").Value2 = 10
$ws.Range("B469").Value2 = "Vega Modelo de Temuco"
$ws.Range("C469").Value2 = "La Araucanía"
$ws.Range("D469").Value2 = 44505
$ws.Range("E469").Value2 = 9
$ws.Range("F469").Value2 = "Fruta"
$ws.Range("G469").Value2 = 100102
$ws.Range("H469").Value2 = "Cítricos"
$ws.Range("I469").Value2 = 100102005
$ws.Range("J469").Value2 = "Naranja"
$ws.Range("K469").Value2 = "Navel Late"
$ws.Range("L469").Value2 = "Primera"
$ws.Range("M469").Value2 = 235
$ws.Range("N469").Value2 = 9000
$ws.Range("O469").Value2 = 10000
$ws.Range("P469").Value2 = 9532
$ws.Range("Q469").Value2 = "`$/caja 15 kilos granel"
$ws.Range("R469").Value2 = "Región de O'Higgins"
$ws.Range("S469").Value2 = 635
$ws.Range("T469").Value2 = 15

# New row 470: Navel Late / Segunda, $/caja 15 kilos granel
$ws.Range("A470").Value2 = 10
$ws.Range("B470").Value2 = "Vega Modelo de Temuco"
$ws.Range("C470").Value2 = "La Araucanía"
$ws.Range("D470").Value2 = 44505
$ws.Range("E470").Value2 = 9
$ws.Range("F470").Value2 = "Fruta"
$ws.Range("G470").Value2 = 100102
$ws.Range("H470").Value2 = "Cítricos"
$ws.Range("I470").Value2 = 100102005
$ws.Range("J470").Value2 = "Naranja"
$ws.Range("K470").Value2 = "Navel Late"
$ws.Range("L470").Value2 = "Segunda"
$ws.Range("M470").Value2 = 155
$ws.Range("N470").Value2 = 8000
$ws.Range("O470").Value2 = 8000
$ws.Range("P470").Value2 = 8000
$ws.Range("Q470").Value2 = "`$/caja 15 kilos granel"
$ws.Range("R470").Value2 = "Región de O'Higgins"
$ws.Range("S470").Value2 = 533
$ws.Range("T470").Value2 = 15
